# Landmark data update: rows 2-113 (landmark indices 0-111) in columns A (index), B (x), C (y).
# New dimension becomes A1:D113 (was A1:D79); column D has no data beyond the header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# x,y pairs for landmark rows 0..111 (row r = landmark index (r-2), columns B,C)
$xy = @(
    @(301.580696105957,198.6799335479736),
    @(304.3289184570312,197.5370121002197),
    @(304.7513389587402,197.4017429351807),
    @(305.1544570922852,197.4861001968384),
    @(303.6117172241211,197.1533918380737),
    @(305.9846687316895,198.3648204803467),
    @(306.4811706542969,199.904580116272),
    @(307.416934967041,202.4381446838379),
    @(310.5128479003906,204.5040035247803),
    @(310.7219886779785,205.3522253036499),
    @(311.2709808349609,205.8853912353516),
    @(311.3628387451172,207.3226976394653),
    @(312.2221946716309,206.9467449188232),
    @(312.3532104492188,208.0878353118896),
    @(312.8471565246582,208.3736371994019),
    @(314.3346977233887,209.8037910461426),
    @(316.6686630249023,209.0056085586548),
    @(317.9023742675781,210.2333736419678),
    @(324.2375564575195,212.1127080917358),
    @(324.1802215576172,210.7591581344604),
    @(314.9545478820801,201.7688941955566),
    @(295.0443267822266,191.7860269546509),
    @(274.5041656494141,185.2526092529297),
    @(253.277587890625,185.6079626083374),
    @(235.5479049682617,209.4093990325928),
    @(221.2269020080566,208.8348627090454),
    @(207.7327728271484,210.8592367172241),
    @(201.8914031982422,218.3521127700806),
    @(201.1476707458496,228.155779838562),
    @(205.6510543823242,240.5823612213135),
    @(212.9120063781738,249.8719596862793),
    @(222.310791015625,258.7912845611572),
    @(228.3261108398438,264.5145606994629),
    @(243.9634704589844,273.3255386352539),
    @(295.644359588623,290.8289051055908),
    @(317.238712310791,293.3724117279053),
    @(337.2767639160156,301.2234020233154),
    @(336.4676284790039,383.2715892791748),
    @(344.6207427978516,389.3365287780762),
    @(362.9837799072266,307.8590297698975),
    @(347.883186340332,348.1959915161133),
    @(368.4585189819336,306.1020755767822),
    @(367.7604293823242,306.2249851226807),
    @(364.7008895874023,307.5014305114746),
    @(364.1349792480469,309.1973876953125),
    @(361.223258972168,306.9017887115479),
    @(357.4379730224609,306.8339538574219),
    @(357.2548294067383,305.3340339660645),
    @(338.4324645996094,336.7652893066406),
    @(356.0394287109375,304.2288780212402),
    @(356.7615509033203,303.0307960510254),
    @(320.1678848266602,393.3853912353516),
    @(356.8764495849609,303.7191581726074),
    @(356.8891906738281,302.2235584259033),
    @(355.1534652709961,304.0714931488037),
    @(354.8737716674805,303.988151550293),
    @(355.010871887207,304.6842098236084),
    @(353.9480590820312,303.4764862060547),
    @(350.3170394897461,302.4907779693604),
    @(347.1619415283203,300.0793075561523),
    @(340.5611419677734,295.01051902771),
    @(333.5683822631836,293.0414199829102),
    @(321.9160461425781,288.7927722930908),
    @(307.404670715332,286.391429901123),
    @(293.4063529968262,285.0151920318604),
    @(277.347297668457,286.5023517608643),
    @(254.8833465576172,299.1382026672363),
    @(240.1508712768555,309.7138595581055),
    @(231.8825912475586,311.5981006622314),
    @(221.7786407470703,319.5260810852051),
    @(213.0893516540527,328.7232112884521),
    @(209.624137878418,341.1438274383545),
    @(256.7561149597168,493.2778358459473),
    @(332.5244903564453,478.2175540924072),
    @(362.150993347168,480.5669403076172),
    @(365.5131149291992,482.5569534301758),
    @(401.2774276733398,392.2531414031982),
    @(406.2438583374023,388.3743667602539),
    @(409.8681640625,384.2517185211182),
    @(413.0849838256836,383.8685703277588),
    @(416.4176940917969,377.2506237030029),
    @(420.0083160400391,375.6314849853516),
    @(421.7411804199219,371.123514175415),
    @(422.6917266845703,371.2782096862793),
    @(422.4713897705078,371.1333847045898),
    @(420.1962280273438,373.1639671325684),
    @(418.6752319335938,372.9774856567383),
    @(417.9193496704102,373.9857959747314),
    @(418.8997650146484,373.7307357788086),
    @(417.7945709228516,374.33349609375),
    @(417.6186370849609,376.0244464874268),
    @(416.9873809814453,375.0797080993652),
    @(417.0002746582031,373.7006950378418),
    @(418.6467742919922,373.7280750274658),
    @(418.184700012207,373.2118606567383),
    @(418.8584518432617,371.608772277832),
    @(416.9054412841797,368.3235740661621),
    @(412.9043960571289,362.9985237121582),
    @(405.9521102905273,350.3407287597656),
    @(404.2134857177734,343.0961894989014),
    @(396.3188934326172,332.9078578948975),
    @(391.0876846313477,324.2669677734375),
    @(390.6815338134766,318.4580039978027),
    @(384.7146987915039,308.1191253662109),
    @(381.8271255493164,296.0513591766357),
    @(380.3066635131836,289.9652767181396),
    @(375.214729309082,282.5899028778076),
    @(374.1634368896484,292.9549884796143),
    @(343.3008193969727,349.8074054718018),
    @(339.8038101196289,344.6346473693848),
    @(337.8397750854492,360.8515548706055),
    @(334.2683029174805,364.176778793335)
)

for ($i = 0; $i -lt $xy.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $xy[$i][0]
    $ws.Cells.Item($row, 3).Value = $xy[$i][1]
}

Write-Host "Updated landmarks data: rows 2-113 (dimension now A1:D113)"
